# Refresh the cryptos list snapshot (prices / 1h volume deltas), matching the
# GitHub Actions scheduled update. A couple of rows also swap rank position
# (Dai/OKB and the Kaspa/PEPE/Stacks/Fetch.AI block), so Coin + Link are
# rewritten for those rows too.
#
# Numeric-looking Price values are prefixed with a leading apostrophe so
# Excel keeps them as literal text (e.g. "1.00") instead of coercing them to
# a number (which would drop the trailing zero / alter the stored value).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.520.39'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '3.350.29'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''189.52'
$ws.Range('E5').Value = '  +4.23%  '
$ws.Range('D6').Value = '''559.32'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = '3.343.48'
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('E9').Value = '  -1.33%  '
$ws.Range('E10').Value = '  -3.31%  '
$ws.Range('D11').Value = '''0.586'
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('D12').Value = '''47.10'
$ws.Range('E12').Value = '  -1.58%  '
$ws.Range('E13').Value = '  +1.69%  '
$ws.Range('D14').Value = '''8.71'
$ws.Range('E14').Value = '  +1.47%  '
$ws.Range('D15').Value = '3.884.46'
$ws.Range('E15').Value = '  +1.15%  '
$ws.Range('D16').Value = '''603.72'
$ws.Range('E16').Value = '  -5.76%  '
$ws.Range('D17').Value = '66.586.03'
$ws.Range('E17').Value = '  +0.65%  '
$ws.Range('D18').Value = '''18.04'
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('D19').Value = '3.353.47'
$ws.Range('E19').Value = '  +1.42%  '
$ws.Range('E20').Value = '  +1.06%  '
$ws.Range('D21').Value = '''11.08'
$ws.Range('E21').Value = '  -3.76%  '
$ws.Range('D22').Value = '''0.905'
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').Value = '''18.40'
$ws.Range('E23').Value = '  +4.30%  '
$ws.Range('D24').Value = '''5.06'
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').Value = '''100.66'
$ws.Range('E25').Value = '  -6.21%  '
$ws.Range('D26').Value = '''4.01'
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('D27').Value = '''6.08'
$ws.Range('E27').Value = '  +1.05%  '
$ws.Range('D28').Value = '''2.76'
$ws.Range('E28').Value = '  +2.27%  '
$ws.Range('D29').Value = '''9.62'
$ws.Range('E29').Value = '  +0.42%  '
$ws.Range('D30').Value = '''8.73'
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('D31').Value = '''30.79'
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('D32').Value = '''6.72'
$ws.Range('E32').Value = '  +5.14%  '
$ws.Range('D33').Value = '''3.95'
$ws.Range('E33').Value = '  -2.60%  '
$ws.Range('D34').Value = '''590.53'
$ws.Range('E34').Value = '  +7.45%  '
$ws.Range('D35').Value = '''11.07'
$ws.Range('E35').Value = '  -0.38%  '
$ws.Range('E36').Value = '  -0.21%  '
$ws.Range('D37').Value = '3.738.28'
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').Value = '''1.00'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '''56.78'
$ws.Range('E39').Value = '  -1.09%  '
$ws.Range('D40').Value = '''3.53'
$ws.Range('E40').Value = '  +4.31%  '
$ws.Range('D41').Value = '''34.00'
$ws.Range('E41').Value = '  +5.16%  '
$ws.Range('B42').Value = 'PEPE'
$ws.Range('C42').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D42').Value = '0.0₃0712'
$ws.Range('E42').Value = '  -1.07%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '''3.26'
$ws.Range('E43').Value = '  -7.71%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = '''2.70'
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('B45').Value = 'Kaspa'
$ws.Range('C45').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D45').Value = '''0.128'
$ws.Range('E45').Value = '  +0.59%  '
$ws.Range('D46').Value = '''0.343'
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('D47').Value = '''3.37'
$ws.Range('E47').Value = '  +3.58%  '
$ws.Range('D48').Value = '''0.0422'
$ws.Range('E48').Value = '  +1.43%  '
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('D50').Value = '''2.59'
$ws.Range('E50').Value = '  -1.83%  '
$ws.Range('E51').Value = '  +0.30%  '
